$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '68.544.94'
Set-TextValue "E2" '  -0.86%  '
Set-TextValue "D3" '3.901.18'
Set-TextValue "E3" '  +2.75%  '
Set-TextValue "E4" '  -0.03%  '
Set-TextValue "D5" '602.32'
Set-TextValue "E5" '  +0.18%  '
Set-TextValue "D6" '166.05'
Set-TextValue "E6" '  +1.40%  '
Set-TextValue "D7" '3.896.73'
Set-TextValue "E7" '  +2.63%  '
Set-TextValue "E8" '  +0.13%  '
Set-TextValue "E9" '  -1.12%  '
Set-TextValue "E10" '  -1.63%  '
Set-TextValue "D11" '6.40'
Set-TextValue "E11" '  +1.27%  '
Set-TextValue "D12" '0.460'
Set-TextValue "E12" '  -0.08%  '
Set-TextValue "D13" '0.0000254'
Set-TextValue "E13" '  +3.66%  '
Set-TextValue "D14" '37.26'
Set-TextValue "E14" '  +0.01%  '
Set-TextValue "D15" '4.555.83'
Set-TextValue "E15" '  +2.80%  '
Set-TextValue "D16" '3.919.07'
Set-TextValue "E16" '  +3.48%  '
Set-TextValue "D17" '68.662.34'
Set-TextValue "E17" '  -0.83%  '
Set-TextValue "D18" '7.47'
Set-TextValue "E18" '  +1.00%  '
Set-TextValue "D19" '17.14'
Set-TextValue "E19" '  -0.92%  '
Set-TextValue "E20" '  -2.18%  '
Set-TextValue "E21" '  -2.68%  '
Set-TextValue "D22" '486.95'
Set-TextValue "E22" '  -0.35%  '
Set-TextValue "E23" '  +0.32%  '
Set-TextValue "E24" '  +9.99%  '
Set-TextValue "D25" '84.50'
Set-TextValue "E25" '  -0.10%  '
Set-TextValue "E26" '  -1.05%  '
Set-TextValue "E27" '  -1.24%  '
Set-TextValue "D28" '10.11'
Set-TextValue "E28" '  +0.53%  '
Set-TextValue "D30" '2.94'
Set-TextValue "E30" '  -0.91%  '
Set-TextValue "D31" '4.054.40'
Set-TextValue "E31" '  +2.70%  '
Set-TextValue "E32" '  -0.43%  '
Set-TextValue "D33" '7.73'
Set-TextValue "E33" '  -3.96%  '
Set-TextValue "D34" '31.80'
Set-TextValue "E34" '  -0.08%  '
Set-TextValue "D35" '3.855.72'
Set-TextValue "E35" '  +2.99%  '
Set-TextValue "E36" '  -0.26%  '
Set-TextValue "E37" '  +2.12%  '
Set-TextValue "E38" '  +0.28%  '
Set-TextValue "D39" '0.138'
Set-TextValue "E39" '  -1.89%  '
Set-TextValue "E40" '  +5.28%  '
Set-TextValue "D41" '1.00'
Set-TextValue "E41" '  -0.05%  '
Set-TextValue "E42" '  -2.22%  '
Set-TextValue "D43" '430.23'
Set-TextValue "E43" '  +2.03%  '
Set-TextValue "D44" '48.33'
Set-TextValue "E44" '  -0.50%  '
Set-TextValue "E45" '  -0.40%  '
Set-TextValue "D46" '8.52'
Set-TextValue "E46" '  +1.76%  '
Set-TextValue "E47" '  -0.02%  '
Set-TextValue "D48" '142.63'
Set-TextValue "E48" '  +1.44%  '
Set-TextValue "D49" '26.18'
Set-TextValue "E49" '  +8.04%  '
Set-TextValue "D50" '2.806.55'
Set-TextValue "E50" '  -0.72%  '
Set-TextValue "E51" '  +0.56%  '
